$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add row 19 string first so shared-string table gets "0289" before "0283"
$ws.Range("I19").Value = "EKA_METALS_PATCH_0289"

# Add row 18
$ws.Range("B18").Value = 40893
$ws.Range("B18").NumberFormat = $ws.Range("B17").NumberFormat
$ws.Range("I18").Value = "EKA_METALS_PATCH_0283"

# Finish row 19
$ws.Range("B19").Value = 40897
$ws.Range("B19").NumberFormat = $ws.Range("B17").NumberFormat

# Update selection to match new active cell
$ws.Range("I18").Select()
